# Applies the "Interannual update / Major update" restructuring described in
# the commit message: renames the "High Priority break-up" sheet to
# "Interannual update - High Pri" with new (interannual) summary data, and
# adds a new sheet "Major update - High Priority " that keeps the previous
# "High Priority break-up" figures. Also refreshes several summary numbers on
# the other sheets to reflect the newer (memory-optimized) pipeline run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Trends Status
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")

$wsTrends.Range("B2").Value = 3
$wsTrends.Range("C2").Value = 1
$wsTrends.Range("D2").Value = 7
$wsTrends.Range("E2").Value = 1.5

$wsTrends.Range("B3").Value = 6
$wsTrends.Range("C3").Value = 7
$wsTrends.Range("D3").Value = 14
$wsTrends.Range("E3").Value = 10.8

$wsTrends.Range("B4").Value = 26
$wsTrends.Range("C4").Value = 54
$wsTrends.Range("D4").Value = 60.5
$wsTrends.Range("E4").Value = 83.09999999999999

$wsTrends.Range("B5").Value = 5
$wsTrends.Range("C5").Value = 2
$wsTrends.Range("D5").Value = 11.6
$wsTrends.Range("E5").Value = 3.1

$wsTrends.Range("B6").Value = 3
$wsTrends.Range("C6").Value = 1
$wsTrends.Range("D6").Value = 7
$wsTrends.Range("E6").Value = 1.5

$wsTrends.Range("B7").Value = 158
$wsTrends.Range("C7").Value = 217

$wsTrends.Range("B8").Value = 272
$wsTrends.Range("C8").Value = 191

# ---------------------------------------------------------------------------
# Sheet: Priority Status
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")

$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# Sheet: Species qualification
# ---------------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")

$wsQual.Range("A2").Value = "SoIB Assessment"
$wsQual.Range("B3").Value = 201
$wsQual.Range("C3").Value = 43
$wsQual.Range("B4").Value = 282
$wsQual.Range("C4").Value = 65

# ---------------------------------------------------------------------------
# Sheet: High Priority break-up -> split into two sheets
# ---------------------------------------------------------------------------
$wsOldBreakup = $wb.Worksheets.Item("High Priority break-up")

# First add the new "Major update" sheet (placed right after the existing
# one) and copy over the original "High Priority break-up" figures verbatim.
$wsMajor = $wb.Worksheets.Add($null, $wsOldBreakup)
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"
$wsMajor.Range("A1:E1").Font.Bold = $true
$wsMajor.Range("A1:E1").HorizontalAlignment = -4108

$wsMajor.Range("A2").Value = "Trend New"
$wsMajor.Range("B2").Value = 12
$wsMajor.Range("C2").Value = 21.4
$wsMajor.Range("D2").Value = 12
$wsMajor.Range("E2").Value = 21.4

$wsMajor.Range("A3").Value = "IUCN"
$wsMajor.Range("B3").Value = 44
$wsMajor.Range("C3").Value = 78.59999999999999
$wsMajor.Range("D3").Value = 44
$wsMajor.Range("E3").Value = 78.59999999999999

# Rename the original sheet and overwrite its data with the new interannual
# update figures.
$wsOldBreakup.Name = "Interannual update - High Pri"

$wsOldBreakup.Range("A2").Value = "Trend New"
$wsOldBreakup.Range("B2").Value = 61
$wsOldBreakup.Range("C2").Value = 59.2
$wsOldBreakup.Range("D2").Value = 61
$wsOldBreakup.Range("E2").Value = 93.8

$wsOldBreakup.Range("A3").Value = "Trend Different"
$wsOldBreakup.Range("B3").Value = 1
$wsOldBreakup.Range("C3").Value = 1
$wsOldBreakup.Range("D3").ClearContents()
$wsOldBreakup.Range("E3").ClearContents()

$wsOldBreakup.Range("A4").Value = "IUCN"
$wsOldBreakup.Range("B4").Value = 41
$wsOldBreakup.Range("C4").Value = 39.8
$wsOldBreakup.Range("D4").Value = 4
$wsOldBreakup.Range("E4").Value = 6.2
